# Add a new column "Maximum Br field in coil (T)" after the
# "Stored energy density (MJ/m^3)" column (old column C), shifting the
# existing Maximum hoop stress / Current density / Maximum-Minimum radial
# stress / Possible conductor material columns one slot to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D; this shifts old D:G -> E:H and copies the
# formatting of the column immediately to the left (old column C).
$ws.Columns("D").Insert()

# New header for the inserted column.
$ws.Range("D1").Value = "Maximum Br field in coil (T)"

# New Br-field values for the inserted column, keyed by row number.
$brValues = @{
    2  = 3.7
    3  = 9.1999999999999993
    4  = 7.9
    5  = 11.2
    7  = -3.1
    8  = -5.0999999999999996
    9  = -13.4
    11 = -7.6
    12 = -8.6999999999999993
    14 = -10.4
    15 = 14
    17 = -12.7
    18 = -11.6
    19 = 11.4
    21 = -9.5
    22 = -10.5
    23 = 12.4
    25 = -10.199999999999999
    26 = -14.9
    27 = 12.6
    29 = -9
    30 = -9.6999999999999993
    31 = 12.4
    33 = -9.1999999999999993
    34 = -10.8
    35 = 11.5
}

foreach ($row in $brValues.Keys) {
    $ws.Cells.Item($row, 4).Value = $brValues[$row]
}

# Updated "Possible conductor material" text (now column H) per row.
$materialValues = @{
    2  = "REBCO/Nb3Sn"
    3  = "REBCO/Nb3Sn"
    4  = "REBCO/Nb3Sn"
    5  = "REBCO/Nb3Sn"
    7  = "REBCO/Nb3Sn"
    8  = "REBCO/Nb3Sn"
    9  = "REBCO/Nb3Sn"
    11 = "REBCO/Nb3Sn"
    12 = "REBCO/Nb3Sn"
    14 = "REBCO/Nb3Sn"
    15 = "REBCO/Nb3Sn"
    17 = "REBCO/Nb3Sn"
    18 = "REBCO/Nb3Sn"
    19 = "REBCO/Nb3Sn"
    21 = "REBCO/Nb3Sn"
    22 = "REBCO/Nb3Sn"
    23 = "REBCO/Nb3Sn"
    25 = "REBCO/(Nb3Sn?)"
    26 = "REBCO/Nb3Sn"
    27 = "REBCO/Nb3Sn"
    29 = "REBCO"
    30 = "REBCO/Nb3Sn"
    31 = "REBCO/Nb3Sn"
    33 = "REBCO"
    34 = "REBCO/Nb3Sn"
    35 = "REBCO/Nb3Sn"
}

foreach ($row in $materialValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $materialValues[$row]
}

# Match the author's final view/selection state.
$ws.Range("G33").Select()
$ws.Application.ActiveWindow.ScrollRow = 16
